# "Generate Report for Archive"
# The localization status changes from "Ready for handoff" to "In Translation"
# for both tracked files/languages. The Status/zh-cn/de-de columns that hold
# this text are narrowed to fit the new (shorter) value.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 13.4101848602295
$overview.Columns.Item(6).ColumnWidth = 13.4101848602295

# --- zh-cn sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 13.4101848602295

# --- de-de sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 13.4101848602295
